# 新增 2022-Q4 工作表数据 ("feat: add 2022-Q4 data")
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert the new "2022-Q4" sheet right after "总计" (i.e. right
#    before "2022-Q3"), mirroring the existing quarterly sheets.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($q3)
$newSheet.Name = "2022-Q4"

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows (index col A, fund code/name/scale/position/ratio/value, rank)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'001606"
$newSheet.Range("C2").Value = "农银汇理工业4.0灵活配置混合"
$newSheet.Range("D2").Value = "'39.21"
$newSheet.Range("E2").Value = "'81.12"
$newSheet.Range("F2").Value = "'3.23"
$newSheet.Range("G2").Value = "'1.2665"
$newSheet.Range("H2").Value = 7

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'000336"
$newSheet.Range("C3").Value = "农银研究精选混合"
$newSheet.Range("D3").Value = "'38.29"
$newSheet.Range("E3").Value = "'88.29"
$newSheet.Range("F3").Value = "'2.11"
$newSheet.Range("G3").Value = "'0.8079"
$newSheet.Range("H3").Value = 9

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'009686"
$newSheet.Range("C4").Value = "华夏磐利一年定期开放混合A"
$newSheet.Range("D4").Value = "'11.49"
$newSheet.Range("E4").Value = "'92.56"
$newSheet.Range("F4").Value = "'3.65"
$newSheet.Range("G4").Value = "'0.4194"
$newSheet.Range("H4").Value = 10

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'000259"
$newSheet.Range("C5").Value = "农银区间收益混合"
$newSheet.Range("D5").Value = "'4.29"
$newSheet.Range("E5").Value = "'75.37"
$newSheet.Range("F5").Value = "'1.21"
$newSheet.Range("G5").Value = "'0.0519"
$newSheet.Range("H5").Value = 8

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'007138"
$newSheet.Range("C6").Value = "鹏扬元合量化大盘优选股票C"
$newSheet.Range("D6").Value = "'0.59"
$newSheet.Range("E6").Value = "'94.05"
$newSheet.Range("F6").Value = "'3.81"
$newSheet.Range("G6").Value = "'0.0225"
$newSheet.Range("H6").Value = 6

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "'009687"
$newSheet.Range("C7").Value = "华夏磐利一年定期开放混合C"
$newSheet.Range("D7").Value = "'0.46"
$newSheet.Range("E7").Value = "'92.56"
$newSheet.Range("F7").Value = "'3.65"
$newSheet.Range("G7").Value = "'0.0168"
$newSheet.Range("H7").Value = 10

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "'519969"
$newSheet.Range("C8").Value = "长信新利灵活配置混合"
$newSheet.Range("D8").Value = "'0.49"
$newSheet.Range("E8").Value = "'91.04"
$newSheet.Range("F8").Value = "'2.65"
$newSheet.Range("G8").Value = "'0.0130"
$newSheet.Range("H8").Value = 5

$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "'005638"
$newSheet.Range("C9").Value = "农银汇理量化智慧动力混合"
$newSheet.Range("D9").Value = "'0.46"
$newSheet.Range("E9").Value = "'90.00"
$newSheet.Range("F9").Value = "'1.92"
$newSheet.Range("G9").Value = "'0.0088"
$newSheet.Range("H9").Value = 7

$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "'007137"
$newSheet.Range("C10").Value = "鹏扬元合量化大盘优选股票A"
$newSheet.Range("D10").Value = "'0.11"
$newSheet.Range("E10").Value = "'94.05"
$newSheet.Range("F10").Value = "'3.81"
$newSheet.Range("G10").Value = "'0.0042"
$newSheet.Range("H10").Value = 6

# Match formatting of the sibling quarter sheets: bold/bordered header
# row and index column, default formatting for the rest.
$q3ref = $wb.Worksheets.Item("2022-Q3")
$q3ref.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q3ref.Range("A2").Copy()
$newSheet.Range("A2:A10").PasteSpecial(-4122)
$q3ref.Range("B2:H2").Copy()
$newSheet.Range("B2:H10").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert a new row for 2022-Q4
#    at row 2 (shifting the existing quarters down) and renumber the
#    sequential index column.
# ------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Rows.Item(2).Insert()

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 9
$totals.Range("D2").Value = 2.61

for ($r = 3; $r -le 9; $r++) {
    $totals.Cells.Item($r, 1).Value = $r - 2
}

$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)
$totals.Range("B3:D3").Copy()
$totals.Range("B2:D2").PasteSpecial(-4122)

Write-Host "2022-Q4 sheet added and 总计 updated"
